$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the Basic Operating Weight (BOW) for the A321neo (LR variant, EFB) block.
$ws.Range("B23").Value = 110500
